# Update the slide's title placeholder text from "Title" to "Slide_title".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

if ($s.Shapes.HasTitle) {
    $titleShape = $s.Shapes.Title
} else {
    $titleShape = $s.Shapes.Item(4)
}

$titleShape.TextFrame.TextRange.Text = "Slide_title"
